$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B47 to be a real number instead of a text value
$ws.Range("B47").Value = 1

# Add new row 48 with the additional annotation data
$ws.Range("A48").Value = "Ruilin"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "3"
$ws.Range("C48").Value = "无"
$ws.Range("D48").Value = "DIS"
$ws.Range("E48").Value = "OTH"
$ws.Range("F48").Value = "3bf3a8cd-f7a3-492e-815a-c1d9e74634b1"
$ws.Range("G48").Value = "ByCPHrgCW_annotated.xlsx"
$ws.Range("H48").Value = "The problem scenario states that the model/weights is private, but later on it ceases to be so (weights are not encrypted)."
